# Add three new numbered-list paragraphs right after the paragraph that
# ends with "Add bibliography / tools use cases / components. Bookmarks,
# Lectures. Notes: Scrapbook.":
#   1) an empty bullet paragraph
#   2) "CQRS. Monads Functions (domain / range) CUD Commands, R Retrievals
#      applicable in contexts / roles: DDD (signatures / dataflow)."
#   3) "CoSQL. Duals. Meijer. LinQ / DSL / Parser Combinators. Templates.
#      Parse Model (Statements): Sets Model AST CUD, Parsed Model
#      Execution: R. Scala Cats."
#
# The three paragraphs share the same list formatting as their neighbours
# (numId 4 / ilvl 0, ind left=600 hanging=360) and carry the paragraph-mark
# rPr (<w:u w:val="none"/>) that Word stamps on freshly split paragraphs.

$d = $word.ActiveDocument

# Locate the anchor paragraph by its distinctive text. Keep a handle to the
# same Range object throughout - re-deriving $d.Content would reset to the
# top of the document and lose the Find position.
$anchorText = "Add bibliography / tools use cases / components. Bookmarks, Lectures. Notes: Scrapbook."
$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text."
}
$anchor = $searchRange.Paragraphs(1)
$insertionPoint = $d.Range($anchor.Range.End, $anchor.Range.End)

# Build the OOXML fragment for the three new paragraphs. A trailing empty
# <w:p/> is appended as a break terminator so InsertXML doesn't merge the
# third paragraph's run into the pre-existing (empty) paragraph that
# follows the anchor; that terminator paragraph is deleted afterwards.
$pPrCommon = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:ind w:left="600" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr>'

$para1 = '<w:p>' + $pPrCommon + '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$para2 = '<w:p>' + $pPrCommon + '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">CQRS. Monads Functions (domain / range) CUD Commands, R Retrievals applicable in contexts / roles: DDD (signatures / dataflow).</w:t></w:r></w:p>'
$para3 = '<w:p>' + $pPrCommon + '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">CoSQL. Duals. Meijer. LinQ / DSL / Parser Combinators. Templates. Parse Model (Statements): Sets Model AST CUD, Parsed Model Execution: R. Scala Cats.</w:t></w:r></w:p>'
$terminator = '<w:p/>'

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $para1 + $para2 + $para3 + $terminator +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$countBefore = $d.Paragraphs.Count
$insertionPoint.InsertXML($xmlFragment)

# Remove the terminator paragraph InsertXML had to introduce to keep the
# third new paragraph from merging into the next (pre-existing) paragraph.
$anchorIndex = $anchor.Index
$terminatorParagraph = $d.Paragraphs($anchorIndex + 4)
$terminatorParagraph.Range.Delete()

Write-Output ("Paragraphs before: " + $countBefore + ", after: " + $d.Paragraphs.Count)
